# Insert a new row of data at row 524 ("Hortaliza, Vega Modelo de Temuco - Perejil"),
# pushing all existing rows 524-554 down to 525-555, matching a new weekly
# price-update commit. The dimension grows from A1:R554 to A1:R555.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 524; this shifts old rows
# 524..554 down to 525..555 and keeps the inherited "date" style (s="2")
# on the D column of the new row.
$ws.Rows.Item(524).Insert()

# Populate the freshly inserted row 524 with the new weekly record.
$ws.Range("A524").Value = 10
$ws.Range("B524").Value = "Vega Modelo de Temuco"
$ws.Range("C524").Value = "La Araucanía"
$ws.Range("D524").Value = 45265
$ws.Range("E524").Value = 9
$ws.Range("F524").Value = 100112044
$ws.Range("G524").Value = "Perejil"
$ws.Range("H524").Value = "Sin especificar"
$ws.Range("I524").Value = "Primera"
$ws.Range("J524").Value = 20
$ws.Range("K524").Value = 6000
$ws.Range("L524").Value = 6000
$ws.Range("M524").Value = 6000
$ws.Range("N524").Value = "$/docena de atados (3 kilos)"
$ws.Range("O524").Value = "Provincia de Cautín"
$ws.Range("P524").Value = 2000
$ws.Range("Q524").Value = 3
$ws.Range("R524").Value = "Hortaliza"
